# Add a new worksheet "Sheet1" to the workbook, mirroring the structure
# used on "format" but demonstrating the new "list as leaf" feature.

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the existing "format" sheet
$wsFormat = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $wsFormat)
$ws2.Name = "Sheet1"

# Header row
$ws2.Range("A1").Value = "####"
$ws2.Range("C1").Value = "listLeaf#list[0]"
$ws2.Range("D1").Value = "listLeaf#list[1]"

# Data rows (list[0] and list[1] leaf values)
$ws2.Range("A3").Value = "listLeaf#~"
$ws2.Range("C4").Value = "test1-0"
$ws2.Range("D3").Value = "test0-1"
$ws2.Range("C3").Value = "test0-0"
$ws2.Range("D4").Value = "test1-1"

# Column widths to match bestFit widths recorded in the diff
$ws2.Columns.Item(3).ColumnWidth = 14.08203125
$ws2.Columns.Item(4).ColumnWidth = 13.9140625

# Selection on the new sheet, and make it the active/selected tab
$ws2.Range("C5").Select() | Out-Null

# Make the new sheet the active tab (workbookView activeTab=1 / sheet2 tabSelected)
$ws2.Activate() | Out-Null
